# Update the "cryptos" worksheet with freshly scraped Price / Volume(1h)
# figures, plus a ranking swap (Monero <-> LidoDAOToken) for rows 25/26.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the Price column to remain plain text so values like "1.000",
# "314.44" or "27.411.03" are not re-interpreted as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.411.03'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '1.827.94'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '314.44'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4468'
$ws.Range('E7').Value = '  +5.04%  '
$ws.Range('D8').Value = '0.3762'
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('D9').Value = '0.07536'
$ws.Range('E9').Value = '  +4.14%  '
$ws.Range('D10').Value = '0.8928'
$ws.Range('E10').Value = '  +5.91%  '
$ws.Range('D11').Value = '21.07'
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('D12').Value = '1.818.05'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '6.752'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').Value = '94.42'
$ws.Range('E14').Value = '  +5.32%  '
$ws.Range('D15').Value = '5.416'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '0.07116'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '0.000008825'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '15.27'
$ws.Range('E20').Value = '  +2.66%  '
$ws.Range('D21').Value = '27.395.66'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').Value = '5.276'
$ws.Range('D23').Value = '10.93'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').Value = '1.978'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = '2.368'
$ws.Range('E25').Value = '  +6.85%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '151.46'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').Value = '18.59'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').Value = '5.377'
$ws.Range('E28').Value = '  +2.85%  '
$ws.Range('D29').Value = '117.72'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').Value = '0.08848'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').Value = '0.7835'
$ws.Range('E31').Value = '  +6.05%  '
$ws.Range('D32').Value = '1.203'
$ws.Range('E32').Value = '  +2.22%  '
$ws.Range('D33').Value = '4.518'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('D34').Value = '2.893'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').Value = '1.001'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '1.109'
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('D37').Value = '0.01991'
$ws.Range('E37').Value = '  +2.60%  '
$ws.Range('D38').Value = '0.05325'
$ws.Range('E38').Value = '  +2.11%  '
$ws.Range('D39').Value = '7.367'
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('D40').Value = '0.5315'
$ws.Range('E40').Value = '  +3.60%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('D42').Value = '0.1730'
$ws.Range('E42').Value = '  +2.45%  '
$ws.Range('D43').Value = '2.295'
$ws.Range('E43').Value = '  +17.92%  '
$ws.Range('D44').Value = '8.762'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('D45').Value = '0.5151'
$ws.Range('E45').Value = '  +8.84%  '
$ws.Range('D46').Value = '10.72'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('D47').Value = '1.706'
$ws.Range('E47').Value = '  +3.15%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = '1.000'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '0.06377'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '0.9355'
$ws.Range('E51').Value = '  +3.50%  '
